$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Plain text values (not at risk of numeric auto-conversion): direct assignment ---
$ws.Range('D2').Value = '30.137.89'
$ws.Range('E2').Value = '  +0.51%  '
$ws.Range('D3').Value = '1.907.70'
$ws.Range('E3').Value = '  +2.12%  '
$ws.Range('E4').Value = '  -0.25%  '
$ws.Range('E5').Value = '  +0.55%  '
$ws.Range('E6').Value = '  -0.20%  '
$ws.Range('E7').Value = '  -0.46%  '
$ws.Range('E8').Value = '  +4.31%  '
$ws.Range('E9').Value = '  +1.97%  '
$ws.Range('E10').Value = '  +1.00%  '
$ws.Range('E11').Value = '  +1.40%  '
$ws.Range('E12').Value = '  +6.03%  '
$ws.Range('B13').Value = 'WrappedEther'
$ws.Range('C13').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D13').Value = '1.908.22'
$ws.Range('E13').Value = '  +1.83%  '
$ws.Range('B14').Value = 'Polkadot'
$ws.Range('C14').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('E14').Value = '  +2.18%  '
$ws.Range('E15').Value = '  +0.92%  '
$ws.Range('E16').Value = '  -0.30%  '
$ws.Range('E17').Value = '  +0.75%  '
$ws.Range('E18').Value = '  +1.65%  '
$ws.Range('E19').Value = '  +2.73%  '
$ws.Range('E20').Value = '  +3.19%  '
$ws.Range('E21').Value = '  -0.34%  '
$ws.Range('E22').Value = '  +2.59%  '
$ws.Range('D23').Value = '30.144.73'
$ws.Range('E23').Value = '  +0.55%  '
$ws.Range('E24').Value = '  +2.74%  '
$ws.Range('E25').Value = '  -0.65%  '
$ws.Range('D26').Value = '2.124.09'
$ws.Range('E26').Value = '  +1.74%  '
$ws.Range('E27').Value = '  +4.25%  '
$ws.Range('E28').Value = '  +1.02%  '
$ws.Range('E29').Value = '  +2.75%  '
$ws.Range('E31').Value = '  +9.73%  '
$ws.Range('E32').Value = '  +1.04%  '
$ws.Range('E33').Value = '  +1.50%  '
$ws.Range('E34').Value = '  +1.14%  '
$ws.Range('E35').Value = '  +1.53%  '
$ws.Range('E36').Value = '  +3.00%  '
$ws.Range('E37').Value = '  +0.79%  '
$ws.Range('E38').Value = '  +0.87%  '
$ws.Range('E39').Value = '  +4.42%  '
$ws.Range('E40').Value = '  +2.61%  '
$ws.Range('E41').Value = '  +1.33%  '
$ws.Range('E42').Value = '  +1.72%  '
$ws.Range('E43').Value = '  -0.20%  '
$ws.Range('E44').Value = '  +5.39%  '
$ws.Range('E45').Value = '  +10.67%  '
$ws.Range('E46').Value = '  +3.42%  '
$ws.Range('E47').Value = '  -0.42%  '
$ws.Range('E48').Value = '  +0.55%  '
$ws.Range('E49').Value = '  -0.20%  '
$ws.Range('E50').Value = '  +3.18%  '
$ws.Range('E51').Value = '  +1.46%  '

# --- Numeric-looking text values: force text via formula + paste-special (values only),
# then clear the formula remnants, to avoid Excel coercing these into real numbers ---
$ws.Range('D4').Formula = '="1.001"'
$ws.Range('D4').Copy()
$ws.Range('D4').PasteSpecial(-4163)
$ws.Range('D5').Formula = '="319.87"'
$ws.Range('D5').Copy()
$ws.Range('D5').PasteSpecial(-4163)
$ws.Range('D6').Formula = '="1.000"'
$ws.Range('D6').Copy()
$ws.Range('D6').PasteSpecial(-4163)
$ws.Range('D8').Formula = '="0.4083"'
$ws.Range('D8').Copy()
$ws.Range('D8').PasteSpecial(-4163)
$ws.Range('D9').Formula = '="0.08333"'
$ws.Range('D9').Copy()
$ws.Range('D9').PasteSpecial(-4163)
$ws.Range('D10').Formula = '="42.35"'
$ws.Range('D10').Copy()
$ws.Range('D10').PasteSpecial(-4163)
$ws.Range('D12').Formula = '="23.84"'
$ws.Range('D12').Copy()
$ws.Range('D12').PasteSpecial(-4163)
$ws.Range('D14').Formula = '="6.387"'
$ws.Range('D14').Copy()
$ws.Range('D14').PasteSpecial(-4163)
$ws.Range('D15').Formula = '="7.223"'
$ws.Range('D15').Copy()
$ws.Range('D15').PasteSpecial(-4163)
$ws.Range('D16').Formula = '="1.000"'
$ws.Range('D16').Copy()
$ws.Range('D16').PasteSpecial(-4163)
$ws.Range('D17').Formula = '="92.37"'
$ws.Range('D17').Copy()
$ws.Range('D17').PasteSpecial(-4163)
$ws.Range('D18').Formula = '="0.00001096"'
$ws.Range('D18').Copy()
$ws.Range('D18').PasteSpecial(-4163)
$ws.Range('D19').Formula = '="0.06493"'
$ws.Range('D19').Copy()
$ws.Range('D19').PasteSpecial(-4163)
$ws.Range('D20').Formula = '="18.38"'
$ws.Range('D20').Copy()
$ws.Range('D20').PasteSpecial(-4163)
$ws.Range('D21').Formula = '="0.9996"'
$ws.Range('D21').Copy()
$ws.Range('D21').PasteSpecial(-4163)
$ws.Range('D22').Formula = '="5.928"'
$ws.Range('D22').Copy()
$ws.Range('D22').PasteSpecial(-4163)
$ws.Range('D24').Formula = '="11.34"'
$ws.Range('D24').Copy()
$ws.Range('D24').PasteSpecial(-4163)
$ws.Range('D25').Formula = '="2.191"'
$ws.Range('D25').Copy()
$ws.Range('D25').PasteSpecial(-4163)
$ws.Range('D28').Formula = '="162.72"'
$ws.Range('D28').Copy()
$ws.Range('D28').PasteSpecial(-4163)
$ws.Range('D29').Formula = '="2.284"'
$ws.Range('D29').Copy()
$ws.Range('D29').PasteSpecial(-4163)
$ws.Range('D30').Formula = '="128.43"'
$ws.Range('D30').Copy()
$ws.Range('D30').PasteSpecial(-4163)
$ws.Range('D31').Formula = '="1.140"'
$ws.Range('D31').Copy()
$ws.Range('D31').PasteSpecial(-4163)
$ws.Range('D32').Formula = '="0.1042"'
$ws.Range('D32').Copy()
$ws.Range('D32').PasteSpecial(-4163)
$ws.Range('D33').Formula = '="5.956"'
$ws.Range('D33').Copy()
$ws.Range('D33').PasteSpecial(-4163)
$ws.Range('D35').Formula = '="0.02454"'
$ws.Range('D35').Copy()
$ws.Range('D35').PasteSpecial(-4163)
$ws.Range('D36').Formula = '="5.348"'
$ws.Range('D36').Copy()
$ws.Range('D36').PasteSpecial(-4163)
$ws.Range('D37').Formula = '="0.06375"'
$ws.Range('D37').Copy()
$ws.Range('D37').PasteSpecial(-4163)
$ws.Range('D38').Formula = '="0.2146"'
$ws.Range('D38').Copy()
$ws.Range('D38').PasteSpecial(-4163)
$ws.Range('D39').Formula = '="0.6528"'
$ws.Range('D39').Copy()
$ws.Range('D39').PasteSpecial(-4163)
$ws.Range('D40').Formula = '="1.195"'
$ws.Range('D40').Copy()
$ws.Range('D40').PasteSpecial(-4163)
$ws.Range('D41').Formula = '="8.606"'
$ws.Range('D41').Copy()
$ws.Range('D41').PasteSpecial(-4163)
$ws.Range('D42').Formula = '="11.36"'
$ws.Range('D42').Copy()
$ws.Range('D42').PasteSpecial(-4163)
$ws.Range('D43').Formula = '="1.209"'
$ws.Range('D43').Copy()
$ws.Range('D43').PasteSpecial(-4163)
$ws.Range('D44').Formula = '="13.46"'
$ws.Range('D44').Copy()
$ws.Range('D44').PasteSpecial(-4163)
$ws.Range('D45').Formula = '="2.195"'
$ws.Range('D45').Copy()
$ws.Range('D45').PasteSpecial(-4163)
$ws.Range('D46').Formula = '="0.6070"'
$ws.Range('D46').Copy()
$ws.Range('D46').PasteSpecial(-4163)
$ws.Range('D47').Formula = '="3.618"'
$ws.Range('D47').Copy()
$ws.Range('D47').PasteSpecial(-4163)
$ws.Range('D48').Formula = '="1.208"'
$ws.Range('D48').Copy()
$ws.Range('D48').PasteSpecial(-4163)
$ws.Range('D49').Formula = '="121.54"'
$ws.Range('D49').Copy()
$ws.Range('D49').PasteSpecial(-4163)
$ws.Range('D50').Formula = '="78.95"'
$ws.Range('D50').Copy()
$ws.Range('D50').PasteSpecial(-4163)
$ws.Range('D51').Formula = '="1.140"'
$ws.Range('D51').Copy()
$ws.Range('D51').PasteSpecial(-4163)
$excel.CutCopyMode = 0

